$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.026619
$ws.Range("H2").Value = 0.079857
$ws.Range("Q2").Value = 0.328906298314
$ws.Range("R2").Value = 2.960156684826
